# Prototype_Interaction.pptx edit script
# 1) Update the datetimeFigureOut placeholder text (05/07/2021 -> 05/10/2021)
#    on the slide master and every slide layout.
# 2) Re-flow the legend / connector shapes on slide 1 and widen+shift the
#    "Hi" process box so there is room for an extra legend row, matching
#    the new "User intent" row already present on slide 2. Also slide the
#    "User intent" legend row up on slide 2 to close the resulting gap.

$p = $ppt.ActivePresentation

# ---- Slide master date placeholder ----
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "05/07/2021") {
            $sh.TextFrame.TextRange.Text = "05/10/2021"
        }
    }
}

# ---- Every slide layout's date placeholder ----
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    $layoutShapes = $layout.Shapes
    for ($j = 1; $j -le $layoutShapes.Count; $j++) {
        $sh = $layoutShapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "05/07/2021") {
                $sh.TextFrame.TextRange.Text = "05/10/2021"
            }
        }
    }
}

# ---- Slide 1 geometry updates ----
$s1 = $p.Slides.Item(1)

# "Hi" process box - shift left & widen (right edge stays put)
$hiBox = $s1.Shapes.Item(1)
$hiBox.Left = 713066 / 12700
$hiBox.Width = 2496927 / 12700

# Legend row 1 (User utterance) - connector + textbox shift down
$s1.Shapes.Item(4).Top = 820724 / 12700
$s1.Shapes.Item(5).Top = 820724 / 12700

# Legend row 2 (Agent utterance) - connector + textbox shift down
$s1.Shapes.Item(6).Top = 1275519 / 12700
$s1.Shapes.Item(7).Top = 1275519 / 12700

# Legend row 3 (Custom action taken by agent) - connector + textbox shift down
$s1.Shapes.Item(13).Top = 1730314 / 12700
$s1.Shapes.Item(14).Top = 1738486 / 12700

# ---- Slide 2 geometry updates ----
$s2 = $p.Slides.Item(2)

# Legend row 4 (User intent) - connector + textbox shift up
$s2.Shapes.Item(11).Top = 2235171 / 12700
$s2.Shapes.Item(12).Top = 2235171 / 12700
